$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (word) and Column C (count) updates per row
$ws.Range("B2").Value = "<there>"
$ws.Range("C2").Value = 29

$ws.Range("B4").Value = "<now>"
$ws.Range("C4").Value = 33

$ws.Range("B5").Value = "<que>"
$ws.Range("C5").Value = 31

$ws.Range("B6").Value = "<not>"
$ws.Range("C6").Value = 30

$ws.Range("B7").Value = "<other>"

$ws.Range("B8").Value = "<for>"
$ws.Range("C8").Value = 28

$ws.Range("C9").Value = 34

$ws.Range("B10").Value = "<all>"
$ws.Range("C10").Value = 28

$ws.Range("B11").Value = "<been>"
$ws.Range("C11").Value = 28

$ws.Range("C12").Value = 25

$ws.Range("C13").Value = 32

$ws.Range("B14").Value = "<six>"
$ws.Range("C14").Value = 31

$ws.Range("B15").Value = "<it>"

$ws.Range("C16").Value = 35

$ws.Range("B17").Value = "<sero>"
$ws.Range("C17").Value = 28

$ws.Range("C18").Value = 30
